$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.759.28"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "2.242.82"
$ws.Range("E3").Value = "  -3.26%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'298.42"
$ws.Range("E5").Value = "  -2.42%  "
$ws.Range("D6").Value = "'83.89"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "'0.518"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").Value = "'30.59"
$ws.Range("E10").Value = "  +3.20%  "
$ws.Range("D12").Value = "'47.37"
$ws.Range("E12").Value = "  -9.78%  "
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("D14").Value = "2.585.50"
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("D15").Value = "'6.35"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "'14.30"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "2.236.44"
$ws.Range("E17").Value = "  -3.29%  "
$ws.Range("D18").Value = "'0.724"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").Value = "39.712.27"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "0.0₃0882"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").Value = "'5.82"
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("D22").Value = "'65.35"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("D23").Value = "'10.52"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").Value = "'229.30"
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("E26").Value = "  -3.78%  "
$ws.Range("D27").Value = "'1.85"
$ws.Range("E27").Value = "  +3.98%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  +2.68%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").Value = "'32.76"
$ws.Range("E31").Value = "  -3.60%  "
$ws.Range("D32").Value = "'150.18"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("D35").Value = "'2.43"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").Value = "'16.23"
$ws.Range("E37").Value = "  +6.35%  "
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("D39").Value = "'0.0975"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("D40").Value = "'2.68"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("D41").Value = "'1.68"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").Value = "'3.75"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "1.927.77"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "'0.0266"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").Value = "'2.05"
$ws.Range("E45").Value = "  -9.21%  "
$ws.Range("D46").Value = "'16.78"
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("D47").Value = "'9.16"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").Value = "2.453.98"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("D50").Value = "'71.56"
$ws.Range("E50").Value = "  +2.86%  "
$ws.Range("D51").Value = "'89.36"
$ws.Range("E51").Value = "  -2.70%  "
